$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Decsion Tree" typo in the CARBON ONLY table (row 8).
$ws.Range("A8").Value = "Decision Tree"

# CARBON and SILICON table: apply min-max scaling results for all four
# models (was only kNN / Random Forest Classifier before).
$ws.Range("A20").Value = "Logistic Regression"
$ws.Range("B20").Value = 87.2
$ws.Range("A21").Value = "Decision Tree"
$ws.Range("B21").Value = 93.2
$ws.Range("A22").Value = "kNN"
$ws.Range("B22").Value = 94.8
$ws.Range("A23").Value = "Random Forest Classifier"
$ws.Range("B23").Value = 95.2

# Move the selection to match where the author ended up after typing the
# new rows in.
[void]$ws.Range("C23").Select()
